# edit.ps1 - applies the "Actualizacion Casos de uso" changes:
#  1) Paragraph "[INCLUDE] Ingresar datos." -> split into three runs:
#       "[INCLUDE] Ingresar datos" + " de ciudadano y volumen de materiales" + "."
#  2) The first empty paragraph right after "Caso de uso incluido: Ingresar datos."
#     gets new content: bold "Suposición:" + " El usuario ya está registrado."

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "[INCLUDE] Ingresar datos." -> split into three runs
# ---------------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.ClearFormatting()
$found1 = $find1.Find.Execute("[INCLUDE] Ingresar datos.", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

if (-not ($found1 -eq $true -or $find1.Find.Found)) {
    # Fallback: locate by paragraph text directly
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        if ($para.Range.Text -eq "[INCLUDE] Ingresar datos.`r") {
            $find1 = $d.Range($para.Range.Start, $para.Range.End - 1)
            break
        }
    }
}

$target1 = $d.Range($find1.Start, $find1.End)
$target1.Text = ""

$insertPoint1 = $d.Range($target1.Start, $target1.Start)
$bodyXml1 = '<w:p><w:r><w:t>[INCLUDE] Ingresar datos</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> de ciudadano y volumen de materiales</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r></w:p>'
$pkg1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $bodyXml1 + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$insertPoint1.InsertXML($pkg1)

# ---------------------------------------------------------------------------
# 2) Fill the first empty paragraph right after
#    "Caso de uso incluido: Ingresar datos." with "Suposición: ..."
# ---------------------------------------------------------------------------
$anchorIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Caso de uso incluido: Ingresar datos.`r") {
        $anchorIdx = $i
        break
    }
}

$targetPara = $d.Paragraphs.Item($anchorIdx + 1)
$insertPoint2 = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
$bodyXml2 = '<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Suposición:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> El usuario ya está registrado.</w:t></w:r></w:p>'
$pkg2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $bodyXml2 + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$insertPoint2.InsertXML($pkg2)

Write-Host "Para8=[$($d.Paragraphs.Item(8).Range.Text)]"
Write-Host "Para13=[$($d.Paragraphs.Item(13).Range.Text)]"
